# Apply updated market/profit values to the Leve profit sheets.
# Generated from the authoritative cell-level diff (before -> after).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 883.2963
$ws.Range("I2").Value = 269.66666
$ws.Range("K2").Value = 269.66666
$ws.Range("M2").Value = -156.66666
$ws.Range("H33").Value = 264.75
$ws.Range("I33").Value = 264.75
$ws.Range("K33").Value = 264.75
$ws.Range("M33").Value = -35.75
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").ClearContents()
$ws.Range("H113").Value = 5000
$ws.Range("I113").Value = 5000
$ws.Range("K113").Value = 5000
$ws.Range("M113").Value = -1746
$ws.Range("H137").Value = 1998
$ws.Range("I137").Value = 1997.5
$ws.Range("K137").Value = 5992.5
$ws.Range("M137").Value = -3442.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2260.625
$ws.Range("I2").Value = 2297.8572
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 2297.8572
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -2184.8572
$ws.Range("N2").Value = -2226
$ws.Range("H26").Value = 1066.3334
$ws.Range("I26").Value = 1066.3334
$ws.Range("K26").Value = 1066.3334
$ws.Range("M26").Value = -736.3334
$ws.Range("H32").Value = 5137789.5
$ws.Range("I32").Value = 5389619
$ws.Range("J32").Value = 3500897.5
$ws.Range("K32").Value = 5389619
$ws.Range("L32").Value = 3500897.5
$ws.Range("M32").Value = -5389332
$ws.Range("N32").Value = -3501471.5
$ws.Range("H61").Value = 2202.1428
$ws.Range("I61").Value = 2166.8333
$ws.Range("K61").Value = 2166.8333
$ws.Range("M61").Value = -1954.8333
$ws.Range("H74").Value = 5224.5
$ws.Range("I74").Value = 5224.5
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 5224.5
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -4350.5
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 5224.5
$ws.Range("I77").Value = 5224.5
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 26122.5
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -21754.5
$ws.Range("N77").ClearContents()
$ws.Range("H110").Value = 7400958
$ws.Range("I110").Value = 7400958
$ws.Range("K110").Value = 7400958
$ws.Range("M110").Value = -7398913
$ws.Range("H116").Value = 2260.625
$ws.Range("I116").Value = 2297.8572
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 2297.8572
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = -3.857199999999921
$ws.Range("N116").Value = -6588
$ws.Range("H122").Value = 1294.3334
$ws.Range("I122").Value = 1019.93335
$ws.Range("K122").Value = 3059.80005
$ws.Range("M122").Value = -609.8000499999998
$ws.Range("H136").Value = 2202.1428
$ws.Range("I136").Value = 2166.8333
$ws.Range("K136").Value = 6500.499899999999
$ws.Range("M136").Value = -3950.499899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2260.625
$ws.Range("I3").Value = 2297.8572
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 2297.8572
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -2183.8572
$ws.Range("N3").Value = -2228
$ws.Range("H22").Value = 155.6
$ws.Range("I22").Value = 155.6
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 155.6
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 17.40000000000001
$ws.Range("N22").ClearContents()
$ws.Range("H105").Value = 2358.3333
$ws.Range("I105").Value = 2600
$ws.Range("J105").Value = 2116.6667
$ws.Range("K105").Value = 2600
$ws.Range("L105").Value = 2116.6667
$ws.Range("M105").Value = -853
$ws.Range("N105").Value = -5610.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 373.25
$ws.Range("I22").Value = 373.25
$ws.Range("K22").Value = 373.25
$ws.Range("M22").Value = -23.25
$ws.Range("H58").Value = 2229
$ws.Range("I58").Value = 1953
$ws.Range("K58").Value = 1953
$ws.Range("M58").Value = -1750
$ws.Range("H86").Value = 14682.25
$ws.Range("I86").Value = 13743
$ws.Range("K86").Value = 13743
$ws.Range("M86").Value = -12620
$ws.Range("H89").Value = 14682.25
$ws.Range("I89").Value = 13743
$ws.Range("K89").Value = 68715
$ws.Range("M89").Value = -63099
$ws.Range("H99").Value = 3000
$ws.Range("I99").Value = 3000
$ws.Range("K99").Value = 3000
$ws.Range("M99").Value = -1502
$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530
$ws.Range("H136").Value = 2229
$ws.Range("I136").Value = 1953
$ws.Range("K136").Value = 5859
$ws.Range("M136").Value = -3309
$ws.Range("H141").Value = 37352.08
$ws.Range("J141").Value = 37352.08
$ws.Range("L141").Value = 37352.08
$ws.Range("N141").Value = -47712.08

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 157.1875
$ws.Range("I12").Value = 197.125
$ws.Range("J12").Value = 117.25
$ws.Range("K12").Value = 591.375
$ws.Range("L12").Value = 351.75
$ws.Range("M12").Value = -418.375
$ws.Range("N12").Value = -697.75
$ws.Range("H92").Value = 499.75
$ws.Range("J92").Value = 599.5
$ws.Range("L92").Value = 1798.5
$ws.Range("N92").Value = -4294.5
$ws.Range("H122").Value = 378.1111
$ws.Range("I122").Value = 233.16667
$ws.Range("J122").Value = 668
$ws.Range("K122").Value = 2098.50003
$ws.Range("L122").Value = 6012
$ws.Range("M122").Value = 351.4999699999998
$ws.Range("N122").Value = -10912
$ws.Range("H131").Value = 2578.9
$ws.Range("I131").Value = 995
$ws.Range("J131").Value = 2974.875
$ws.Range("K131").Value = 2985
$ws.Range("L131").Value = 8924.625
$ws.Range("M131").Value = 2055
$ws.Range("N131").Value = -19004.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 12468.5
$ws.Range("J18").Value = 20037
$ws.Range("L18").Value = 20037
$ws.Range("N18").Value = -20381
$ws.Range("H22").Value = 549.1
$ws.Range("I22").Value = 493.25
$ws.Range("K22").Value = 493.25
$ws.Range("M22").Value = -198.25
$ws.Range("H27").Value = 549.1
$ws.Range("I27").Value = 493.25
$ws.Range("K27").Value = 493.25
$ws.Range("M27").Value = -386.25
$ws.Range("H46").Value = 2119.6924
$ws.Range("I46").Value = 895.2222
$ws.Range("K46").Value = 895.2222
$ws.Range("M46").Value = -707.2222
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H122").Value = 8515.559999999999
$ws.Range("I122").Value = 8607.923000000001
$ws.Range("K122").Value = 25823.769
$ws.Range("M122").Value = -23373.769

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2299.1428
$ws.Range("I136").Value = 2415.6667
$ws.Range("J136").Value = 1600
$ws.Range("K136").Value = 7247.000100000001
$ws.Range("L136").Value = 4800
$ws.Range("M136").Value = -4697.000100000001
$ws.Range("N136").Value = -9900

Write-Host "Applied all cell updates."